$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: C54, 1000pF ---
$ws.Range("F23").Value = "445-1308-1-ND"
$ws.Range("G23").Value = "C1608X7R1H102K"
$ws.Range("E23").Value = "DK"
$ws.Range("I23").Value = 0.1

# --- Row 10: C18, 0.01uF ---
$ws.Range("L10").Value = "0.01uF = 10k pF"
$ws.Range("F10").Value = "478-1227-1-ND"
$ws.Range("G10").Value = "06035C103KAT2A"
$ws.Range("E10").Value = "DK"
$ws.Range("I10").Value = 0.06

# --- Row 12: C20, 0.22uF ---
$ws.Range("F12").Value = "445-5191-1-ND"
$ws.Range("G12").Value = "C1608X7R1E224K"
$ws.Range("E12").Value = "DK"
$ws.Range("I12").Value = 0.28

# --- Row 13: C26, 10nF (reuses row 10's vendor/mfr PN) ---
$ws.Range("L13").Value = "10nF = 10k pF (redundant line item)"
$ws.Range("F13").Value = "478-1227-1-ND"
$ws.Range("G13").Value = "06035C103KAT2A"
$ws.Range("E13").Value = "DK"
$ws.Range("I13").Value = 0.06

# --- Row 14: C27, 470nF ---
$ws.Range("L14").Value = "470nF = 0.47uF"
$ws.Range("F14").Value = "445-3454-1-ND"
$ws.Range("G14").Value = "C1608Y5V1E474Z"
$ws.Range("E14").Value = "DK"
$ws.Range("I14").Value = 0.17

# --- Row 15: C31, C43, C46, 33nF ---
$ws.Range("L15").Value = "33nF = 33k pF = 0.033uF"
$ws.Range("F15").Value = "445-5106-1-ND"
$ws.Range("G15").Value = "C1608X7R1E333K"
$ws.Range("E15").Value = "DK"
$ws.Range("I15").Value = 0.16

# Column L width change (target stored width 32.85546875; COM ColumnWidth
# quantizes internally, so 32 is the closest achievable input)
$ws.Range("L1").EntireColumn.ColumnWidth = 32

# Selection change
$ws.Range("E16").Select()
